$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 51.22717533333333
$ws.Range("H2").Value = 153.681526
$ws.Range("I2").Value = 0.1770805608477903
$ws.Range("J2").Value = 0.1770805608477904
$ws.Range("M2").Value = 9.020393666666665
$ws.Range("N2").Value = 27.061181
$ws.Range("O2").Value = 0.2369163967891402
$ws.Range("P2").Value = 0.2369163967891402
$ws.Range("Q2").Value = 462.0892879380228
$ws.Range("R2").Value = 4158.803591442205
$ws.Range("S2").Value = 0.04195328841745859
$ws.Range("T2").Value = 0.04195328841745859

$ws.Range("G3").Value = 51.22717533333333
$ws.Range("H3").Value = 153.681526
$ws.Range("I3").Value = 0.1770805608477903
$ws.Range("J3").Value = 0.1770805608477904
$ws.Range("M3").Value = 6.554651666666667
$ws.Range("N3").Value = 19.663955
$ws.Range("O3").Value = 0.1721548429547032
$ws.Range("P3").Value = 0.1721548429547032
$ws.Range("Q3").Value = 335.7762901772589
$ws.Range("R3").Value = 3021.98661159533
$ws.Range("S3").Value = 0.03048527614308211
$ws.Range("T3").Value = 0.03048527614308212

$ws.Range("G4").Value = 51.22717533333333
$ws.Range("H4").Value = 153.681526
$ws.Range("I4").Value = 0.1770805608477903
$ws.Range("J4").Value = 0.1770805608477904
$ws.Range("M4").Value = 10.017966
$ws.Range("N4").Value = 30.053898
$ws.Range("O4").Value = 0.2631171649023133
$ws.Range("P4").Value = 0.2631171649023133
$ws.Range("Q4").Value = 513.1921007653719
$ws.Range("R4").Value = 4618.728906888347
$ws.Range("S4").Value = 0.04659293512958218
$ws.Range("T4").Value = 0.04659293512958219

$ws.Range("G5").Value = 51.22717533333333
$ws.Range("H5").Value = 153.681526
$ws.Range("I5").Value = 0.1770805608477903
$ws.Range("J5").Value = 0.1770805608477904
$ws.Range("M5").Value = 1.524802333333333
$ws.Range("N5").Value = 4.574407
$ws.Range("O5").Value = 0.04004821607331257
$ws.Range("P5").Value = 0.04004821607331257
$ws.Range("Q5").Value = 78.11131647834245
$ws.Range("R5").Value = 703.0018483050819
$ws.Range("S5").Value = 0.007091760563215681
$ws.Range("T5").Value = 0.007091760563215682

$ws.Range("G6").Value = 51.22717533333333
$ws.Range("H6").Value = 153.681526
$ws.Range("I6").Value = 0.1770805608477903
$ws.Range("J6").Value = 0.1770805608477904
$ws.Range("M6").Value = 10.95635
$ws.Range("N6").Value = 32.86905
$ws.Range("O6").Value = 0.2877633792805307
$ws.Range("P6").Value = 0.2877633792805307
$ws.Range("Q6").Value = 561.2628624633667
$ws.Range("R6").Value = 5051.3657621703
$ws.Range("S6").Value = 0.05095730059445178
$ws.Range("T6").Value = 0.05095730059445179

$ws.Range("G7").Value = 77.026568
$ws.Range("H7").Value = 231.079704
$ws.Range("I7").Value = 0.2662631264141754
$ws.Range("J7").Value = 0.2662631264141754
$ws.Range("M7").Value = 9.020393666666665
$ws.Range("N7").Value = 27.061181
$ws.Range("O7").Value = 0.2369163967891402
$ws.Range("P7").Value = 0.2369163967891402
$ws.Range("Q7").Value = 694.8099661522692
$ws.Range("R7").Value = 6253.289695370423
$ws.Range("S7").Value = 0.06308210050785777
$ws.Range("T7").Value = 0.06308210050785779

$ws.Range("G8").Value = 77.026568
$ws.Range("H8").Value = 231.079704
$ws.Range("I8").Value = 0.2662631264141754
$ws.Range("J8").Value = 0.2662631264141754
$ws.Range("M8").Value = 6.554651666666667
$ws.Range("N8").Value = 19.663955
$ws.Range("O8").Value = 0.1721548429547032
$ws.Range("P8").Value = 0.1721548429547032
$ws.Range("Q8").Value = 504.8823223188134
$ws.Range("R8").Value = 4543.94090086932
$ws.Range("S8").Value = 0.04583848671246066
$ws.Range("T8").Value = 0.04583848671246066

$ws.Range("G9").Value = 77.026568
$ws.Range("H9").Value = 231.079704
$ws.Range("I9").Value = 0.2662631264141754
$ws.Range("J9").Value = 0.2662631264141754
$ws.Range("M9").Value = 10.017966
$ws.Range("N9").Value = 30.053898
$ws.Range("O9").Value = 0.2631171649023133
$ws.Range("P9").Value = 0.2631171649023133
$ws.Range("Q9").Value = 771.6495393206879
$ws.Range("R9").Value = 6944.845853886191
$ws.Range("S9").Value = 0.07005839894012408
$ws.Range("T9").Value = 0.0700583989401241

$ws.Range("G10").Value = 77.026568
$ws.Range("H10").Value = 231.079704
$ws.Range("I10").Value = 0.2662631264141754
$ws.Range("J10").Value = 0.2662631264141754
$ws.Range("M10").Value = 1.524802333333333
$ws.Range("N10").Value = 4.574407
$ws.Range("O10").Value = 0.04004821607331257
$ws.Range("P10").Value = 0.04004821607331257
$ws.Range("Q10").Value = 117.4502906150587
$ws.Range("R10").Value = 1057.052615535528
$ws.Range("S10").Value = 0.01066336321899064
$ws.Range("T10").Value = 0.01066336321899064

$ws.Range("G11").Value = 77.026568
$ws.Range("H11").Value = 231.079704
$ws.Range("I11").Value = 0.2662631264141754
$ws.Range("J11").Value = 0.2662631264141754
$ws.Range("M11").Value = 10.95635
$ws.Range("N11").Value = 32.86905
$ws.Range("O11").Value = 0.2877633792805307
$ws.Range("P11").Value = 0.2877633792805307
$ws.Range("Q11").Value = 843.9300383068
$ws.Range("R11").Value = 7595.3703447612
$ws.Range("S11").Value = 0.07662077703474224
$ws.Range("T11").Value = 0.07662077703474225

$ws.Range("G12").Value = 72.76991766666667
$ws.Range("H12").Value = 218.309753
$ws.Range("I12").Value = 0.2515488654100336
$ws.Range("J12").Value = 0.2515488654100336
$ws.Range("M12").Value = 9.020393666666665
$ws.Range("N12").Value = 27.061181
$ws.Range("O12").Value = 0.2369163967891402
$ws.Range("P12").Value = 0.2369163967891402
$ws.Range("Q12").Value = 656.4133044442548
$ws.Range("R12").Value = 5907.719739998292
$ws.Range("S12").Value = 0.05959605080934156
$ws.Range("T12").Value = 0.05959605080934156

$ws.Range("G13").Value = 72.76991766666667
$ws.Range("H13").Value = 218.309753
$ws.Range("I13").Value = 0.2515488654100336
$ws.Range("J13").Value = 0.2515488654100336
$ws.Range("M13").Value = 6.554651666666667
$ws.Range("N13").Value = 19.663955
$ws.Range("O13").Value = 0.1721548429547032
$ws.Range("P13").Value = 0.1721548429547032
$ws.Range("Q13").Value = 476.9814621170129
$ws.Range("R13").Value = 4292.833159053115
$ws.Range("S13").Value = 0.04330535542009812
$ws.Range("T13").Value = 0.04330535542009812

$ws.Range("G14").Value = 72.76991766666667
$ws.Range("H14").Value = 218.309753
$ws.Range("I14").Value = 0.2515488654100336
$ws.Range("J14").Value = 0.2515488654100336
$ws.Range("M14").Value = 10.017966
$ws.Range("N14").Value = 30.053898
$ws.Range("O14").Value = 0.2631171649023133
$ws.Range("P14").Value = 0.2631171649023133
$ws.Range("Q14").Value = 729.006561007466
$ws.Range("R14").Value = 6561.059049067194
$ws.Range("S14").Value = 0.06618682430108165
$ws.Range("T14").Value = 0.06618682430108165

$ws.Range("G15").Value = 72.76991766666667
$ws.Range("H15").Value = 218.309753
$ws.Range("I15").Value = 0.2515488654100336
$ws.Range("J15").Value = 0.2515488654100336
$ws.Range("M15").Value = 1.524802333333333
$ws.Range("N15").Value = 4.574407
$ws.Range("O15").Value = 0.04004821607331257
$ws.Range("P15").Value = 0.04004821607331257
$ws.Range("Q15").Value = 110.9597402546079
$ws.Range("R15").Value = 998.637662291471
$ws.Range("S15").Value = 0.01007408331493765
$ws.Range("T15").Value = 0.01007408331493765

$ws.Range("G16").Value = 72.76991766666667
$ws.Range("H16").Value = 218.309753
$ws.Range("I16").Value = 0.2515488654100336
$ws.Range("J16").Value = 0.2515488654100336
$ws.Range("M16").Value = 10.95635
$ws.Range("N16").Value = 32.86905
$ws.Range("O16").Value = 0.2877633792805307
$ws.Range("P16").Value = 0.2877633792805307
$ws.Range("Q16").Value = 797.2926874271834
$ws.Range("R16").Value = 7175.634186844651
$ws.Range("S16").Value = 0.07238655156457467
$ws.Range("T16").Value = 0.07238655156457467

$ws.Range("G17").Value = 32.02005133333333
$ws.Range("H17").Value = 96.06015400000001
$ws.Range("I17").Value = 0.110685951579145
$ws.Range("J17").Value = 0.110685951579145
$ws.Range("M17").Value = 9.020393666666665
$ws.Range("N17").Value = 27.061181
$ws.Range("O17").Value = 0.2369163967891402
$ws.Range("P17").Value = 0.2369163967891402
$ws.Range("Q17").Value = 288.8334682535415
$ws.Range("R17").Value = 2599.501214281874
$ws.Range("S17").Value = 0.02622331682330828
$ws.Range("T17").Value = 0.02622331682330828

$ws.Range("G18").Value = 32.02005133333333
$ws.Range("H18").Value = 96.06015400000001
$ws.Range("I18").Value = 0.110685951579145
$ws.Range("J18").Value = 0.110685951579145
$ws.Range("M18").Value = 6.554651666666667
$ws.Range("N18").Value = 19.663955
$ws.Range("O18").Value = 0.1721548429547032
$ws.Range("P18").Value = 0.1721548429547032
$ws.Range("Q18").Value = 209.8802828387856
$ws.Range("R18").Value = 1888.92254554907
$ws.Range("S18").Value = 0.01905512261139959
$ws.Range("T18").Value = 0.01905512261139959

$ws.Range("G19").Value = 32.02005133333333
$ws.Range("H19").Value = 96.06015400000001
$ws.Range("I19").Value = 0.110685951579145
$ws.Range("J19").Value = 0.110685951579145
$ws.Range("M19").Value = 10.017966
$ws.Range("N19").Value = 30.053898
$ws.Range("O19").Value = 0.2631171649023133
$ws.Range("P19").Value = 0.2631171649023133
$ws.Range("Q19").Value = 320.775785575588
$ws.Range("R19").Value = 2886.982070180292
$ws.Range("S19").Value = 0.02912337377401937
$ws.Range("T19").Value = 0.02912337377401937

$ws.Range("G20").Value = 32.02005133333333
$ws.Range("H20").Value = 96.06015400000001
$ws.Range("I20").Value = 0.110685951579145
$ws.Range("J20").Value = 0.110685951579145
$ws.Range("M20").Value = 1.524802333333333
$ws.Range("N20").Value = 4.574407
$ws.Range("O20").Value = 0.04004821607331257
$ws.Range("P20").Value = 0.04004821607331257
$ws.Range("Q20").Value = 48.82424898651978
$ws.Range("R20").Value = 439.418240878678
$ws.Range("S20").Value = 0.004432774905121811
$ws.Range("T20").Value = 0.004432774905121812

$ws.Range("G21").Value = 32.02005133333333
$ws.Range("H21").Value = 96.06015400000001
$ws.Range("I21").Value = 0.110685951579145
$ws.Range("J21").Value = 0.110685951579145
$ws.Range("M21").Value = 10.95635
$ws.Range("N21").Value = 32.86905
$ws.Range("O21").Value = 0.2877633792805307
$ws.Range("P21").Value = 0.2877633792805307
$ws.Range("Q21").Value = 350.8228894259667
$ws.Range("R21").Value = 3157.406004833701
$ws.Range("S21").Value = 0.03185136346529596
$ws.Range("T21").Value = 0.03185136346529596

$ws.Range("G22").Value = 56.24368933333333
$ws.Range("H22").Value = 168.731068
$ws.Range("I22").Value = 0.1944214957488556
$ws.Range("J22").Value = 0.1944214957488557
$ws.Range("M22").Value = 9.020393666666665
$ws.Range("N22").Value = 27.061181
$ws.Range("O22").Value = 0.2369163967891402
$ws.Range("P22").Value = 0.2369163967891402
$ws.Range("Q22").Value = 507.3402190523674
$ws.Range("R22").Value = 4566.061971471308
$ws.Range("S22").Value = 0.04606164023117403
$ws.Range("T22").Value = 0.04606164023117403

$ws.Range("G23").Value = 56.24368933333333
$ws.Range("H23").Value = 168.731068
$ws.Range("I23").Value = 0.1944214957488556
$ws.Range("J23").Value = 0.1944214957488557
$ws.Range("M23").Value = 6.554651666666667
$ws.Range("N23").Value = 19.663955
$ws.Range("O23").Value = 0.1721548429547032
$ws.Range("P23").Value = 0.1721548429547032
$ws.Range("Q23").Value = 368.6577920282156
$ws.Range("R23").Value = 3317.92012825394
$ws.Range("S23").Value = 0.03347060206766275
$ws.Range("T23").Value = 0.03347060206766275

$ws.Range("G24").Value = 56.24368933333333
$ws.Range("H24").Value = 168.731068
$ws.Range("I24").Value = 0.1944214957488556
$ws.Range("J24").Value = 0.1944214957488557
$ws.Range("M24").Value = 10.017966
$ws.Range("N24").Value = 30.053898
$ws.Range("O24").Value = 0.2631171649023133
$ws.Range("P24").Value = 0.2631171649023133
$ws.Range("Q24").Value = 563.4473674558959
$ws.Range("R24").Value = 5071.026307103064
$ws.Range("S24").Value = 0.05115563275750606
$ws.Range("T24").Value = 0.05115563275750607

$ws.Range("G25").Value = 56.24368933333333
$ws.Range("H25").Value = 168.731068
$ws.Range("I25").Value = 0.1944214957488556
$ws.Range("J25").Value = 0.1944214957488557
$ws.Range("M25").Value = 1.524802333333333
$ws.Range("N25").Value = 4.574407
$ws.Range("O25").Value = 0.04004821607331257
$ws.Range("P25").Value = 0.04004821607331257
$ws.Range("Q25").Value = 85.76050873074178
$ws.Range("R25").Value = 771.844578576676
$ws.Range("S25").Value = 0.007786234071046792
$ws.Range("T25").Value = 0.007786234071046793

$ws.Range("G26").Value = 56.24368933333333
$ws.Range("H26").Value = 168.731068
$ws.Range("I26").Value = 0.1944214957488556
$ws.Range("J26").Value = 0.1944214957488557
$ws.Range("M26").Value = 10.95635
$ws.Range("N26").Value = 32.86905
$ws.Range("O26").Value = 0.2877633792805307
$ws.Range("P26").Value = 0.2877633792805307
$ws.Range("Q26").Value = 616.2255456272667
$ws.Range("R26").Value = 5546.0299106454
$ws.Range("S26").Value = 0.05594738662146603
$ws.Range("T26").Value = 0.05594738662146604
